$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 4 (pushes existing rows 4-7 down to 5-8),
# copying formatting from the row above (row 3) as Excel does by default.
$ws.Rows.Item(4).Insert()

# Fill in the new compound row: 2,4,5-trichlorophenol (MW 197.4)
$ws.Range("A4").Value2 = "2,4,5-trichlorophenol"
$ws.Range("B4").Value2 = 197.4
$ws.Range("C4").Value2 = 5
$ws.Range("D4").Value2 = 10
$ws.Range("E4").Value2 = 20
$ws.Range("F4").Value2 = 50
$ws.Range("G4").Value2 = "n.a."
$ws.Range("H4").Value2 = "n.a."
$ws.Range("I4").Value2 = 155710
$ws.Range("J4").Value2 = 343277
$ws.Range("K4").Value2 = 805095
$ws.Range("L4").Value2 = 2302730
$ws.Range("M4").Value2 = "n.a."
$ws.Range("N4").Value2 = "n.a."
